# Applies odds/score updates for rows 2, 3, 6, 7, 8 as per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes
$ws.Range("G2").Value = 4.2
$ws.Range("H2").Value = 3.7
$ws.Range("I2").Value = 1.77
$ws.Range("J2").Value = 4.75
$ws.Range("L2").Value = 2.4
$ws.Range("AA2").Value = 15

# Row 3 changes
$ws.Range("G3").Value = 5.5
$ws.Range("H3").Value = 3.75
$ws.Range("I3").Value = 1.62
$ws.Range("K3").Value = 2.25
$ws.Range("L3").Value = 2.2
$ws.Range("S3").Value = 3.25
$ws.Range("T3").Value = 1.33
$ws.Range("U3").Value = 1.36
$ws.Range("V3").Value = 3
$ws.Range("Z3").Value = 29
$ws.Range("AE3").Value = 11
$ws.Range("AF3").Value = 7.5
$ws.Range("AG3").Value = 15
$ws.Range("AJ3").Value = 8

# Row 6 changes
$ws.Range("G6").Value = 1.07
$ws.Range("J6").Value = 1.33
$ws.Range("L6").Value = 17
$ws.Range("W6").Value = 2.2
$ws.Range("X6").Value = 1.62
$ws.Range("AC6").Value = 12
$ws.Range("AD6").Value = 34
$ws.Range("AE6").Value = 26
$ws.Range("AF6").Value = 21
$ws.Range("AJ6").Value = 126
$ws.Range("AL6").Value = 451
$ws.Range("AO6").Value = 501

# Row 7 changes
$ws.Range("G7").Value = 2.35
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 2.88
$ws.Range("J7").Value = 2.88
$ws.Range("L7").Value = 3.4
$ws.Range("M7").Value = 1.04
$ws.Range("N7").Value = 13
$ws.Range("O7").Value = 1.2
$ws.Range("P7").Value = 4.33
$ws.Range("Q7").Value = 1.7
$ws.Range("R7").Value = 2.1
$ws.Range("S7").Value = 2.63
$ws.Range("T7").Value = 1.44
$ws.Range("U7").Value = 1.33
$ws.Range("V7").Value = 3.25
$ws.Range("W7").Value = 1.57
$ws.Range("X7").Value = 2.25
$ws.Range("Y7").Value = 10
$ws.Range("Z7").Value = 13
$ws.Range("AA7").Value = 9.5
$ws.Range("AE7").Value = 13
$ws.Range("AF7").Value = 6.5
$ws.Range("AG7").Value = 12
$ws.Range("AH7").Value = 41
$ws.Range("AJ7").Value = 17
$ws.Range("AK7").Value = 11
$ws.Range("AM7").Value = 21
$ws.Range("AN7").Value = 26
$ws.Range("AO7").Value = 126
$ws.Range("AP7").Value = 2
$ws.Range("AQ7").Value = 1.85

# Row 8 changes
$ws.Range("G8").Value = 1.38
$ws.Range("H8").Value = 5
$ws.Range("I8").Value = 7
$ws.Range("J8").Value = 1.83
$ws.Range("K8").Value = 2.5
$ws.Range("O8").Value = 1.18
$ws.Range("P8").Value = 4.5
$ws.Range("Q8").Value = 1.6
$ws.Range("R8").Value = 2.3
$ws.Range("S8").Value = 2.5
$ws.Range("T8").Value = 1.5
$ws.Range("Y8").Value = 8
$ws.Range("Z8").Value = 7
$ws.Range("AB8").Value = 9.5
$ws.Range("AE8").Value = 15
$ws.Range("AF8").Value = 9.5
$ws.Range("AG8").Value = 19
$ws.Range("AI8").Value = 19
$ws.Range("AP8").Value = 2.03
$ws.Range("AQ8").Value = 1.83
